$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 5-14: column A/C formulas change denominator from (n-1) to (n+1),
#     and columns B/D get refreshed simulation-error values ---
$ws.Range("A5").Formula = "=(2*PI())/(10+1)"
$ws.Range("C5").Formula = "=1/(10+1)"
$ws.Range("B5").Value = 0.1207
$ws.Range("D5").Value = 0.00000000000000022204

$ws.Range("A6").Formula = "=(2*PI())/(25+1)"
$ws.Range("C6").Formula = "=1/(25+1)"
$ws.Range("B6").Value = 0.0228
$ws.Range("D6").Value = 0.00000000000000044409

$ws.Range("A7").Formula = "=(2*PI())/(50+1)"
$ws.Range("C7").Formula = "=1/(50+1)"
$ws.Range("B7").Value = 0.006
$ws.Range("D7").Value = 0.0000000000000019984

$ws.Range("A8").Formula = "=(2*PI())/(75+1)"
$ws.Range("C8").Formula = "=1/(75+1)"
$ws.Range("B8").Value = 0.0027
$ws.Range("D8").Value = 0.0000000000000032196

$ws.Range("A9").Formula = "=(2*PI())/(100+1)"
$ws.Range("C9").Formula = "=1/(100+1)"
$ws.Range("B9").Value = 0.0015
$ws.Range("D9").Value = 0.0000000000000017486

$ws.Range("A10").Formula = "=(2*PI())/(150+1)"
$ws.Range("C10").Formula = "=1/(150+1)"
$ws.Range("B10").Value = 0.00068151
$ws.Range("D10").Value = 0.0000000000000047184

$ws.Range("A11").Formula = "=(2*PI())/(250+1)"
$ws.Range("C11").Formula = "=1/(250+1)"
$ws.Range("B11").Value = 0.00024661
$ws.Range("D11").Value = 0.0000000000000032196

$ws.Range("A12").Formula = "=(2*PI())/(500+1)"
$ws.Range("C12").Formula = "=1/(500+1)"
$ws.Range("B12").Value = 0.000061894
$ws.Range("D12").Value = 0.000000000000025979

$ws.Range("A13").Formula = "=(2*PI())/(750+1)"
$ws.Range("C13").Formula = "=1/(750+1)"
$ws.Range("B13").Value = 0.000027545
$ws.Range("D13").Value = 0.000000000000019651

$ws.Range("A14").Formula = "=(2*PI())/(1000+1)"
$ws.Range("C14").Formula = "=1/(1000+1)"
$ws.Range("B14").Value = 0.000015504
$ws.Range("D14").Value = 0.00000000000014932

# --- The very small error values (column D, and the tail of column B)
#     now display in scientific notation ---
$ws.Range("D5:D14").NumberFormat = "0.00E+00"
$ws.Range("B10:B14").NumberFormat = "0.00E+00"

# --- Selection / zoom as left by the author on re-save ---
[void]$ws.Range("F12").Select()
[void]($excel.ActiveWindow.Zoom = 100)
